$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3943215398501536
$ws.Range("C2").Value = 0.0527517350114266
$ws.Range("E2").Value = 0.4118363392397839
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.8821953957540103
$ws.Range("H2").Value = 0.9047491368009872
$ws.Range("K2").Value = 0.3704063069657764
$ws.Range("B3").Value = 0.3543074924380676
$ws.Range("C3").Value = 0.0468462589372507
$ws.Range("E3").Value = 0.3592817994173032
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.8708663079471251
$ws.Range("H3").Value = 0.9047311778077329
$ws.Range("K3").Value = 0.3271482789718334
$ws.Range("B4").Value = 0.3298725098463819
$ws.Range("C4").Value = 0.04319963404847726
$ws.Range("E4").Value = 0.3271148072235803
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.8645514126151568
$ws.Range("H4").Value = 0.9051962245440137
$ws.Range("K4").Value = 0.3006581513086815
$ws.Range("B5").Value = 0.319948619015662
$ws.Range("C5").Value = 0.04170836536081879
$ws.Range("E5").Value = 0.3140300387691326
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.8621380369472149
$ws.Range("H5").Value = 0.9055049782747204
$ws.Range("K5").Value = 0.2898807807689252
$ws.Range("B6").Value = 0.3183027894471593
$ws.Range("C6").Value = 0.04146042347328205
$ws.Range("E6").Value = 0.3118586857476089
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.8617469247968188
$ws.Range("H6").Value = 0.9055634349114712
$ws.Range("K6").Value = 0.288092263765293
$ws.Range("B7").Value = 0.3297385369917265
$ws.Range("C7").Value = 0.04317954355491338
$ws.Range("E7").Value = 0.3269382487970063
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.8645182188225959
$ws.Range("H7").Value = 0.9051999063240714
$ws.Range("K7").Value = 0.3005127329634263
$ws.Range("B8").Value = 0.3804968594102149
$ws.Range("C8").Value = 0.05071977835839903
$ws.Range("E8").Value = 0.3936931812373956
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.878155289722244
$ws.Range("H8").Value = 0.904643850678255
$ws.Range("K8").Value = 0.3554762044113602
$ws.Range("B9").Value = 0.4811050247502919
$ws.Range("C9").Value = 0.0653453898107017
$ws.Range("E9").Value = 0.5255045689635836
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.9100435444185564
$ws.Range("H9").Value = 0.9073529698789855
$ws.Range("K9").Value = 0.4638314700142416
$ws.Range("B10").Value = 0.5556965594243195
$ws.Range("C10").Value = 0.0759980490102663
$ws.Range("E10").Value = 0.6230434594760226
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.9366929368805188
$ws.Range("H10").Value = 0.9116916547433505
$ws.Range("K10").Value = 0.5438154295572417
$ws.Range("B11").Value = 0.5897820227657746
$ws.Range("C11").Value = 0.08082529592175547
$ws.Range("E11").Value = 0.6675999972294306
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.9495342093098884
$ws.Range("H11").Value = 0.914182370661166
$ws.Range("K11").Value = 0.5802901490773138
$ws.Range("B12").Value = 0.6027116342969805
$ws.Range("C12").Value = 0.08265063937962225
$ws.Range("E12").Value = 0.684501553253611
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.9545015902254477
$ws.Range("H12").Value = 0.9152004401582019
$ws.Range("K12").Value = 0.5941154213950881
$ws.Range("B13").Value = 0.5999260222168061
$ws.Range("C13").Value = 0.08225763502056793
$ws.Range("E13").Value = 0.6808601761435114
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.9534270986802937
$ws.Range("H13").Value = 0.9149778412774481
$ws.Range("K13").Value = 0.5911373162169866
$ws.Range("B14").Value = 0.590845304269294
$ws.Range("C14").Value = 0.08097552039330935
$ws.Range("E14").Value = 0.6689899048411121
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.9499407730892528
$ws.Range("H14").Value = 0.9142646238393013
$ws.Range("K14").Value = 0.581427299056827
$ws.Range("B15").Value = 0.5852859956216605
$ws.Range("C15").Value = 0.08018984742766122
$ws.Range("E15").Value = 0.6617228606311869
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.9478189700155895
$ws.Range("H15").Value = 0.9138375265225989
$ws.Range("K15").Value = 0.5754813437653752
$ws.Range("B16").Value = 0.5534720916795095
$ws.Range("C16").Value = 0.07568220575016937
$ws.Range("E16").Value = 0.6201355133469519
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.9358683015151712
$ws.Range("H16").Value = 0.9115393307049544
$ws.Range("K16").Value = 0.5414335410331148
$ws.Range("B17").Value = 0.5339946686565611
$ws.Range("C17").Value = 0.07291216504223996
$ws.Range("E17").Value = 0.5946721680093816
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.9287219108923637
$ws.Range("H17").Value = 0.9102622690151918
$ws.Range("K17").Value = 0.520569408402821
$ws.Range("B18").Value = 0.5228062004927381
$ws.Range("C18").Value = 0.07131714697558778
$ws.Range("E18").Value = 0.5800436577191022
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.9246790260483806
$ws.Range("H18").Value = 0.9095763692740206
$ws.Range("K18").Value = 0.5085773333107113
$ws.Range("B19").Value = 0.5190204521528585
$ws.Range("C19").Value = 0.07077679654878466
$ws.Range("E19").Value = 0.5750936036487815
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.9233217312769142
$ws.Range("H19").Value = 0.9093524716151364
$ws.Range("K19").Value = 0.5045184633392239
$ws.Range("B20").Value = 0.5360665792064481
$ws.Range("C20").Value = 0.07320722273964009
$ws.Range("E20").Value = 0.5973809743451568
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.9294756573561642
$ws.Range("H20").Value = 0.9103931772584417
$ws.Range("K20").Value = 0.522789557856953
$ws.Range("B21").Value = 0.5935119280463255
$ws.Range("C21").Value = 0.08135217943897999
$ws.Range("E21").Value = 0.6724756869732715
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.9509619397249196
$ws.Range("H21").Value = 0.9144720762744782
$ws.Range("K21").Value = 0.5842790096088777
$ws.Range("B22").Value = 0.6311851705483775
$ws.Range("C22").Value = 0.08666007779198992
$ws.Range("E22").Value = 0.7217243137598643
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.9656152060644843
$ws.Range("H22").Value = 0.9175746398032345
$ws.Range("K22").Value = 0.6245423083961725
$ws.Range("B23").Value = 0.6110663680203743
$ws.Range("C23").Value = 0.08382853255230316
$ws.Range("E23").Value = 0.6954230906676457
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.9577381432309267
$ws.Range("H23").Value = 0.915878599012899
$ws.Range("K23").Value = 0.6030459687497967
$ws.Range("B24").Value = 0.5351298392971557
$ws.Range("C24").Value = 0.07307383489960273
$ws.Range("E24").Value = 0.5961562898473147
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.9291346841264385
$ws.Range("H24").Value = 0.9103338432273631
$ws.Range("K24").Value = 0.5217858189849096
$ws.Range("B25").Value = 0.4537709686093194
$ws.Range("C25").Value = 0.06140545293828836
$ws.Range("E25").Value = 0.4897342810580199
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.90085735671218
$ws.Range("H25").Value = 0.9062100200256111
$ws.Range("K25").Value = 0.434454553461336
